# Rewrites the "Culinary Crossroads" essay into "The Marvelous World of
# Biology", per the commit's content swap (title/author/email/body text).
$d = $word.ActiveDocument

# Find + Delete + InsertAfter (instead of Find's own Replacement arg) so
# straight apostrophes in the new text are not auto-corrected into curly
# "smart quotes" by the Find/Replace AutoFormat pipeline. The explicit
# Font.* reassignment afterwards guards against the insertion point losing
# its run formatting (rFonts/color/size) when the deleted range reached
# all the way to the end of a paragraph.
function Replace-Text($old, $new, $sizePt) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false)
    if (-not $found) {
        throw "Find failed for: $old"
    }
    $rng.Delete()
    $rng.InsertAfter($new)

    # Re-find the just-inserted text with a fresh Range so the formatting
    # fix below only touches the new run (reusing $rng directly can end up
    # addressing a stale/zero-length position and bleeding into the
    # paragraph mark's rPr).
    $fixRng = $d.Content
    $null = $fixRng.Find.Execute($new, $true, $false, $false, $false, $false, $true, 1, $false)
    $fixRng.Font.Name = "Calibri"
    $fixRng.Font.Color = 0
    if ($sizePt) {
        $fixRng.Font.Size = $sizePt
    }
}

# --- Title / byline / email -------------------------------------------------
Replace-Text "Culinary Crossroads: Exploring Global Connections in Cuisine" "The Marvelous World of Biology: Exploring the Secrets of Life" 22
Replace-Text "Amelia Gomez, PhD" "Alexandria Hayes" 18
Replace-Text "agomez@culinaryinstitute" "alexandriahayes@emailworld" 16
Replace-Text "edu" "net" 16

# --- Intro paragraph (first block of body text) ------------------------------
Replace-Text "As humans, our connection to food transcends mere sustenance; it weaves a rich tapestry of culture, history, and identity" "Biology, the study of life, embarks on an extraordinary journey into the intricate workings of living organisms" 12

Replace-Text " Throughout history, cuisine has served as a conduit for global exchange, bridging diverse regions and fostering intercultural dialogue" " As we delve into the diverse tapestry of nature's marvels, we unravel the mysteries that govern the behavior and interactions of plants, animals, and microorganisms" 12

Replace-Text " From the spice trade that reshaped culinary landscapes to the diaspora of people carrying beloved recipes across borders, food has consistently played a pivotal role in shaping our global community" " From the smallest cell to the vast ecosystems, biology captivates our imagination and challenges us to comprehend the enigmatic dance of life" 12

Replace-Text "This essay delves into the myriad ways in which cuisine has facilitated global interconnectedness, showcasing the fascinating interplay between food, culture, and history" "Biology serves as a bridge between the grandeur of the universe and the intricacies of our own bodies" 12

# This sentence (plus its trailing period) is replaced by two new sentences.
Replace-Text " Whether it be through the dissemination of ingredients, the fusion of culinary techniques, or the influence of migration patterns, food has acted as a potent force in fostering understanding and appreciation among people from vastly different backgrounds." " It unravels the secrets of DNA, the blueprint of heredity, and uncovers the mysteries of genetic variation, the driving force of evolution and adaptation. We witness the enchanting spectacle of cells dividing, organizing, and communicating, their ceaseless symphony orchestrating the wonders of life." 12

Replace-Text "The exploration begins by examining the historical significance of trade routes, particularly the Silk Road, in facilitating the exchange of spices, herbs, and other culinary treasures" "With each passing day, biological discoveries illuminate the path toward medical advancements and environmental solutions" 12

# Likewise replaced by two new sentences.
Replace-Text " This vibrant network of interconnectedness not only introduced novel flavors and ingredients to various regions but also spurred innovation and experimentation within the culinary realm." " We gain invaluable insights into diseases, their causes, and potential treatments, offering hope to those touched by illness. As we explore the interconnections within ecosystems, we unravel the intricate web of life's dependencies and strive for sustainable practices that harmonize human activities with the natural world." 12

# --- Summary paragraph (runs here have no explicit sz -> default size) -------
Replace-Text "This essay unveils the profound role that cuisine has played in fostering global interconnectedness, demonstrating how food has served as a catalyst for cultural exchange, innovation, and understanding among diverse communities" "Biology, the captivating study of life, unveils the complexities and marvels of living organisms" $null

# Replaced by four new sentences.
Replace-Text " The exploration of historical trade routes, the impact of cultural diffusion, and the influence of migration patterns reveals the remarkable ways in which cuisine has transcended geographic boundaries, enriching our collective culinary heritage and fostering a sense of global kinship." " From the enigmatic dance of cells to the majesty of ecosystems, biology inspires awe and wonder. This science uncovers the mysteries of heredity, evolution, and adaptation, illuminating the tapestry of life's astonishing diversity. It plays a pivotal role in medical advancements and environmental solutions, offering hope for a healthier and harmonious world. Biology, a symphony of discovery, invites us to explore the captivating secrets of life, revealing the profound interconnectedness of all living things." $null

# --- Trailing empty paragraph added after the Summary paragraph --------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
